# Add "refrigerator" and "freezer" (plus the existing "refrigerator_freezer")
# choices to the equipment_types list on the "choices" sheet, and move the
# active tab/selection from "survey" to "choices" (cell D18), matching the
# authored commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("choices")

# Insert three new blank rows right after the last existing equipment_types
# row (row 14), before the blank separator row that precedes climate_zones.
# This pushes the old rows 16/17 (climate_zones) down to 19/20, leaving a
# blank row 18 as the new separator - matching the target layout.
$ws.Rows.Item(15).Resize(3).Insert()

# Match the row height used by the other equipment_types rows.
$ws.Rows.Item(15).Resize(3).RowHeight = 20.15

# New choice: refrigerator
$ws.Range("A15").Value = "equipment_types"
$ws.Range("B15").Value = "refrigerator"
$ws.Range("C15").Value = "Refrigerator"
$ws.Range("D15").Value = "Frigorífico"

# New choice: freezer
$ws.Range("A16").Value = "equipment_types"
$ws.Range("B16").Value = "freezer"
$ws.Range("C16").Value = "Freezer"
$ws.Range("D16").Value = "Congelador"

# New choice: refrigerator_freezer
$ws.Range("A17").Value = "equipment_types"
$ws.Range("B17").Value = "refrigerator_freezer"
$ws.Range("C17").Value = "Refrigerator/Freezer"
$ws.Range("D17").Value = "Frigorífico/Congelador"

# Widen column D to fit the longer "Refrigerator/Freezer" style text.
$ws.Columns.Item(4).ColumnWidth = 26.8

# Move the active selection/tab to the choices sheet, cell D18 - this also
# flips the workbook's activeTab / each sheet's tabSelected flag to match.
$ws.Range("D18").Select()
